$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename / extend columns ---
$ws.Range("A1").Value = "order_id"
$ws.Range("B1").Value = "customer_name"
$ws.Range("C1").Value = "cake"
$ws.Range("D1").Value = "quantity"
$ws.Range("E1").Value = "cost_per_cake"
$ws.Range("F1").Value = "sum_of_each_cake"
$ws.Range("G1").Value = "order_date"

# New header cells (E1:G1) should carry the same style as the existing
# header cells (bordered, bold, centered) - copy format from A1.
$ws.Range("A1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)

# --- Data rows ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "rahul"
$ws.Range("C2").Value = "Classic_Chocolate"
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = 300
$ws.Range("G2").Value = 45282.71052083333

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "rahul"
$ws.Range("C3").Value = "Vanilla"
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 800
$ws.Range("G3").Value = 45282.710625

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "raj"
$ws.Range("C4").Value = "Bliss"
$ws.Range("D4").Value = 50
$ws.Range("E4").Value = 35
$ws.Range("F4").Value = 1750
$ws.Range("G4").Value = 45282.71193287037

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "raj"
$ws.Range("C5").Value = "Cookies"
$ws.Range("D5").Value = 15
$ws.Range("E5").Value = 35
$ws.Range("F5").Value = 525
$ws.Range("G5").Value = 45282.7121412037

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "rahul"
$ws.Range("C6").Value = "Cookies"
$ws.Range("D6").Value = 17
$ws.Range("E6").Value = 35
$ws.Range("F6").Value = 595
$ws.Range("G6").Value = 45282.71256944445

# --- Date/time formatting for the new order_date column ---
# Apply to a single cell first: this mirrors how it was done interactively
# in Excel - the first format string ("yyyy-mm-dd h:mm:ss") registers
# numFmtId 164, then re-applying with a different-cased format string
# ("YYYY-MM-DD HH:MM:SS") registers numFmtId 165 and that's the one that
# actually ends up referenced by the cell style.
$ws.Range("G2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("G2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Apply the resulting format to the rest of the date column, one cell at a
# time, so Excel reuses the same cell-format index instead of minting a new
# one for each distinct "before" style combination.
$ws.Range("G3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("G6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
